$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the two new
# header cells before setting their text, so they pick up the same style
# index as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for new columns I (I0) and J (IF), rows 2-36
$data = @(
    @(5, 6),
    @(9, 9),
    @(10, 10),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(10, 10),
    @(8, 8),
    @(7, 7),
    @(6, 8),
    @(6, 7),
    @(7, 8),
    @(5, 7),
    @(1, 5),
    @(1, 5),
    @(6, 9),
    @(7, 9),
    @(7, 9),
    @(7, 7),
    @(5, 5),
    @(1, 3),
    @(10, 11),
    @(5, 7),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(10, 11),
    @(5, 7),
    @(6, 6),
    @(2, 5),
    @(10, 11),
    @(4, 6),
    @(3, 7),
    @(1, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
